$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp header
$ws.Range("A1").Value = "Datos actualizados a 30 de Septiembre de 2020 a las 17:23"

# Full country data table (row, country, casos totales, nuevos casos, casos activos, recuperados, casos criticos, muertes hoy, muertes)
$data = @(
  @(4, "Estados Unidos", 7412124, 5978, 4651017, 2550029, 0, 293, 211078),
  @(5, "India", 6245404, 21885, 5206044, 941599, 0, 232, 97761),
  @(6, "Brasil", 4780317, 0, 4135088, 502219, 0, 0, 143010),
  @(7, "Rusia", 1176286, 8481, 958257, 197307, 0, 177, 20722),
  @(8, "Colombia", 824042, 0, 734154, 64060, 0, 0, 25828),
  @(9, "Peru", 811768, 0, 676925, 102447, 0, 0, 32396),
  @(10, "España", 758172, 0, 0, 0, 0, 0, 31614),
  @(11, "Mexico", 738163, 4446, 530945, 130055, 0, 560, 77163),
  @(12, "Argentina", 736609, 0, 585857, 134233, 0, 0, 16519),
  @(13, "Sudafrica", 672572, 0, 606520, 49385, 0, 0, 16667),
  @(14, "Francia", 550690, 0, 96327, 422470, 0, 0, 31893),
  @(15, "Chile", 462991, 1691, 436589, 13661, 0, 16, 12741),
  @(16, "Iran", 457219, 3582, 380956, 50094, 0, 183, 26169),
  @(17, "Reino Unido", 446156, 0, 0, 0, 0, 0, 42072),
  @(18, "Banglades", 363479, 1436, 275487, 82741, 0, 32, 5251),
  @(19, "Irak", 362981, 4691, 292197, 61603, 0, 59, 9181),
  @(20, "Arabia Saudita", 334605, 418, 319154, 10683, 0, 29, 4768),
  @(21, "Turquia", 317272, 0, 278504, 30638, 0, 0, 8130),
  @(22, "Italia", 314861, 1851, 227704, 51263, 0, 19, 35894),
  @(23, "Pakistan", 312263, 747, 296881, 8903, 0, 5, 6479),
  @(24, "Filipinas", 311694, 2426, 253488, 52702, 0, 58, 5504),
  @(25, "Alemania", 291191, 725, 256000, 25632, 0, 3, 9559),
  @(26, "Indonesia", 287008, 4284, 214947, 61321, 0, 139, 10740),
  @(27, "Israel", 239806, 2880, 173109, 65150, 0, 19, 1547),
  @(28, "Ucrania", 208959, 4027, 92360, 112470, 0, 64, 4129),
  @(29, "Canada", 157586, 625, 134194, 14097, 0, 4, 9295),
  @(30, "Ecuador", 135749, 0, 112296, 12141, 0, 0, 11312),
  @(31, "Bolivia", 134641, 418, 94895, 31815, 0, 31, 7931),
  @(32, "Rumania", 127572, 2158, 102476, 20271, 0, 33, 4825),
  @(33, "Catar", 125760, 227, 122699, 2847, 0, 0, 214),
  @(34, "Marruecos", 121183, 0, 100253, 18778, 0, 0, 2152),
  @(35, "Paises Bajos", 120845, 3294, 0, 0, 0, 13, 6406),
  @(36, "Belgica", 117115, 1762, 19386, 87728, 0, 14, 10001),
  @(37, "Republica Dominicana", 112209, 309, 87629, 22475, 0, 4, 2105),
  @(38, "Panama", 111853, 0, 88202, 21287, 0, 0, 2364),
  @(39, "Kazajistan", 107908, 75, 102874, 3309, 0, 0, 1725),
  @(40, "Kuwait", 105182, 614, 96688, 7884, 0, 3, 610),
  @(41, "Egipto", 103079, 0, 96094, 1071, 0, 0, 5914),
  @(42, "Oman", 98585, 0, 88528, 9122, 0, 0, 935),
  @(43, "Emiratos Arabes Unidos", 93090, 0, 82538, 10136, 0, 0, 416),
  @(44, "Suecia", 92863, 0, 0, 0, 0, 5, 5893),
  @(45, "Guatemala", 91746, 778, 80256, 8244, 0, 8, 3246),
  @(46, "Polonia", 91514, 1552, 69695, 19306, 0, 30, 2513),
  @(47, "China", 85403, 19, 80578, 191, 0, 0, 4634),
  @(48, "Japon", 82494, 0, 75459, 5478, 0, 0, 1557),
  @(49, "Bielorrusia", 78631, 371, 74525, 3273, 0, 5, 833),
  @(50, "Nepal", 77817, 1559, 56428, 20891, 0, 7, 498),
  @(51, "Honduras", 76098, 561, 27383, 46392, 0, 22, 2323),
  @(52, "Portugal", 75542, 825, 48530, 25041, 0, 8, 1971),
  @(53, "Costa Rica", 74604, 0, 30703, 43021, 0, 0, 880),
  @(54, "Etiopia", 74584, 0, 30952, 42441, 0, 0, 1191),
  @(55, "Venezuela", 74363, 0, 64200, 9542, 0, 0, 621),
  @(56, "Barein", 70422, 0, 64267, 5907, 0, 2, 248),
  @(57, "Chequia", 67843, 0, 32759, 34448, 0, 0, 636),
  @(58, "Nigeria", 58647, 0, 49937, 7599, 0, 0, 1111),
  @(59, "Singapur", 57765, 23, 57488, 250, 0, 0, 27),
  @(60, "Uzbekistan", 56519, 165, 52970, 3081, 0, 2, 468),
  @(61, "Suiza", 53282, 411, 42700, 8508, 0, 5, 2074),
  @(62, "Moldavia", 53042, 1013, 38997, 12725, 0, 10, 1320),
  @(63, "Argelia", 51368, 0, 36063, 13579, 0, 0, 1726),
  @(64, "Armenia", 50359, 458, 44001, 5399, 0, 1, 959),
  @(65, "Kirguistan", 46669, 147, 42879, 2726, 0, 0, 1064),
  @(66, "Ghana", 46482, 0, 45651, 530, 0, 0, 301),
  @(67, "Austria", 44813, 772, 35644, 8370, 0, 3, 799),
  @(68, "Azerbaiyan", 40229, 110, 37954, 1684, 0, 1, 591),
  @(69, "Paraguay", 40101, 0, 23748, 15512, 0, 0, 841),
  @(70, "Estado de Palestina", 39899, 358, 31743, 7845, 0, 5, 311),
  @(71, "Afganistan", 39268, 14, 32789, 5021, 0, 0, 1458),
  @(72, "Kenia", 38529, 151, 24908, 12910, 0, 4, 711),
  @(73, "Libano", 38377, 0, 17110, 20906, 0, 0, 361),
  @(74, "Irlanda", 35740, 0, 23364, 10573, 0, 0, 1803),
  @(75, "Libia", 34525, 511, 19361, 14613, 0, 11, 551),
  @(76, "Serbia", 33551, 72, 31536, 1266, 0, 0, 749),
  @(77, "El Salvador", 29077, 96, 23813, 4421, 0, 4, 843),
  @(78, "Dinamarca", 27998, 534, 20754, 6594, 0, 0, 650),
  @(79, "Bosnia y Herzegovina", 27469, 243, 20616, 5997, 0, 13, 856),
  @(80, "Australia", 27078, 15, 24754, 1438, 0, 4, 886),
  @(81, "Hungria", 26461, 894, 5890, 19806, 0, 8, 765),
  @(82, "Corea del Sur", 23812, 113, 21590, 1809, 0, 6, 413),
  @(83, "Camerun", 20838, 0, 19440, 980, 0, 0, 418),
  @(84, "Bulgaria", 20547, 0, 14489, 5245, 0, 0, 813),
  @(85, "Costa de Marfil", 19669, 0, 19241, 308, 0, 0, 120),
  @(86, "Grecia", 18123, 0, 9989, 7746, 0, 0, 388),
  @(87, "Republica de Macedonia", 17977, 191, 14959, 2279, 0, 2, 739),
  @(88, "Tunez", 17405, 0, 5032, 12127, 0, 0, 246),
  @(89, "Croacia", 16593, 213, 15057, 1256, 0, 5, 280),
  @(90, "Madagascar", 16408, 31, 15301, 877, 0, 1, 230),
  @(91, "Senegal", 14982, 37, 12437, 2234, 0, 1, 311),
  @(92, "Zambia", 14759, 44, 13959, 468, 0, 0, 332),
  @(93, "Noruega", 13961, 47, 11190, 2497, 0, 0, 274),
  @(94, "Albania", 13649, 131, 7847, 5415, 0, 3, 387),
  @(95, "Sudan", 13640, 0, 6764, 6040, 0, 0, 836),
  @(96, "Birmania", 13373, 948, 3755, 9308, 0, 26, 310),
  @(97, "Jordania", 11825, 1776, 4626, 7138, 0, 4, 61),
  @(98, "Namibia", 11265, 125, 9014, 2130, 0, 0, 121),
  @(99, "Malasia", 11224, 89, 9967, 1121, 0, 2, 136),
  @(100, "Consejo Danes para los Refugiados", 10659, 28, 10139, 248, 0, 0, 272),
  @(101, "Guinea", 10634, 0, 9960, 608, 0, 0, 66),
  @(102, "Montenegro", 10575, 0, 7002, 3409, 0, 0, 164),
  @(103, "Maldivas", 10194, 0, 9037, 1123, 0, 0, 34),
  @(104, "Eslovaquia", 10141, 567, 4395, 5698, 0, 3, 48),
  @(105, "Finlandia", 9992, 100, 8100, 1548, 0, 0, 344),
  @(106, "Guayana Francesa", 9929, 0, 9569, 294, 0, 0, 66),
  @(107, "Tayikistan", 9769, 43, 8572, 1121, 0, 1, 76),
  @(108, "Gabon", 8752, 0, 7955, 743, 0, 0, 54),
  @(109, "Haiti", 8740, 0, 6757, 1756, 0, 0, 227),
  @(110, "Mozambique", 8556, 0, 5205, 3292, 0, 0, 59),
  @(111, "Luxemburgo", 8431, 0, 7136, 1171, 0, 0, 124),
  @(112, "Uganda", 8129, 112, 4260, 3794, 0, 0, 75),
  @(113, "Zimbabue", 7837, 0, 6122, 1487, 0, 0, 228),
  @(114, "Mauritania", 7488, 0, 7111, 216, 0, 0, 161),
  @(115, "Jamaica", 6482, 74, 1867, 4508, 0, 6, 107),
  @(116, "Georgia", 6192, 326, 3120, 3033, 0, 3, 39),
  @(117, "Cabo Verde", 5900, 0, 5228, 613, 0, 0, 59),
  @(118, "Malaui", 5772, 0, 4245, 1348, 0, 0, 179),
  @(119, "Eslovenia", 5690, 203, 3804, 1736, 0, 1, 150),
  @(120, "Cuba", 5597, 66, 4893, 582, 0, 0, 122),
  @(121, "Suazilandia", 5462, 0, 4859, 495, 0, 0, 108),
  @(122, "Republica de Yibuti", 5416, 0, 5344, 11, 0, 0, 61),
  @(123, "Nicaragua", 5170, 0, 2913, 2106, 0, 0, 151),
  @(124, "Hong Kong", 5088, 8, 4827, 156, 0, 0, 105),
  @(125, "Guinea Ecuatorial", 5030, 0, 4769, 178, 0, 0, 83),
  @(126, "Congo", 5008, 0, 3887, 1032, 0, 0, 89),
  @(127, "Angola", 4905, 0, 1833, 2893, 0, 0, 179),
  @(128, "Surinam", 4863, 0, 4676, 83, 0, 0, 104),
  @(129, "Ruanda", 4836, 0, 3125, 1682, 0, 0, 29),
  @(130, "Republica de Africa Central", 4806, 0, 1840, 2904, 0, 0, 62),
  @(131, "Lituania", 4693, 115, 2365, 2236, 0, 0, 92),
  @(132, "Trinidad yTobago", 4517, 54, 2560, 1882, 0, 1, 75),
  @(133, "Guadalupe", 4487, 0, 2199, 2246, 0, 0, 42),
  @(134, "Siria", 4148, 0, 1088, 2863, 0, 0, 197),
  @(135, "Bahamas", 4022, 119, 2141, 1786, 0, 4, 95),
  @(136, "Reunion", 3993, 111, 2819, 1158, 0, 2, 16),
  @(137, "Aruba", 3934, 0, 3181, 727, 0, 0, 26),
  @(138, "Somalia", 3588, 0, 2946, 543, 0, 0, 99),
  @(139, "Gambia", 3579, 0, 2161, 1306, 0, 0, 112),
  @(140, "Tailandia", 3564, 5, 3374, 131, 0, 0, 59),
  @(141, "Mayotte", 3541, 0, 2964, 537, 0, 0, 40),
  @(142, "Sri Lanka", 3379, 5, 3230, 136, 0, 0, 13),
  @(143, "Estonia", 3371, 57, 2605, 702, 0, 0, 64),
  @(144, "Botsuana", 3172, 0, 710, 2446, 0, 0, 16),
  @(145, "Mali", 3101, 0, 2443, 527, 0, 0, 131),
  @(146, "Malta", 3058, 23, 2562, 462, 0, 0, 34),
  @(147, "Guyana", 2846, 0, 1644, 1124, 0, 0, 78),
  @(148, "Islandia", 2728, 33, 2167, 551, 0, 0, 10),
  @(149, "Sudan del Sur", 2700, 0, 1290, 1361, 0, 0, 49),
  @(150, "Benin", 2357, 17, 1973, 343, 0, 1, 41),
  @(151, "Guinea-Bisau", 2324, 0, 1549, 736, 0, 0, 39),
  @(152, "Sierra Leona", 2222, 0, 1685, 465, 0, 0, 72),
  @(153, "Uruguay", 2033, 0, 1771, 214, 0, 0, 48),
  @(154, "Burkina Faso", 2032, 0, 1309, 665, 0, 0, 58),
  @(155, "Yemen", 2031, 0, 1275, 169, 0, 0, 587),
  @(156, "Principado de Andorra", 1966, 0, 1265, 648, 0, 0, 53),
  @(157, "Belice", 1943, 52, 1225, 692, 0, 2, 26),
  @(158, "Nueva Zelanda", 1836, 1, 1767, 44, 0, 0, 25),
  @(159, "Letonia", 1824, 95, 1307, 480, 0, 0, 37),
  @(160, "Togo", 1759, 0, 1341, 370, 0, 0, 48),
  @(161, "Republica de Chipre", 1743, 0, 1369, 352, 0, 0, 22),
  @(162, "Polinesia Francesa", 1728, 0, 1431, 290, 0, 0, 7),
  @(163, "Lesoto", 1576, 0, 873, 668, 0, 0, 35),
  @(164, "Liberia", 1343, 0, 1221, 40, 0, 0, 82),
  @(165, "Martinica", 1290, 0, 98, 1172, 0, 0, 20),
  @(166, "Niger", 1196, 0, 1110, 17, 0, 0, 69),
  @(167, "Republica del Chad", 1193, 0, 1007, 101, 0, 0, 85),
  @(168, "Vietnam", 1094, 0, 1010, 49, 0, 0, 35),
  @(169, "Santo Tome y Principe", 911, 0, 885, 11, 0, 0, 15),
  @(170, "San Marino", 732, 5, 680, 10, 0, 0, 42),
  @(171, "Crucero", 712, 0, 651, 48, 0, 0, 13),
  @(172, "Islas Turcas y Caicos", 686, 4, 642, 38, 0, 1, 6),
  @(173, "San Martin (Parte Holandesa)", 659, 14, 557, 80, 0, 0, 22),
  @(174, "Papua Nueva Guinea", 534, 0, 516, 11, 0, 0, 7),
  @(175, "Taiwan", 514, 1, 483, 24, 0, 0, 7),
  @(176, "Tanzania", 509, 0, 183, 305, 0, 0, 21),
  @(177, "Burundi", 506, 0, 472, 33, 0, 0, 1),
  @(178, "Comoras", 479, 0, 464, 8, 0, 0, 7),
  @(179, "Islas Feroe", 467, 4, 429, 38, 0, 0, 0),
  @(180, "Gibraltar", 396, 5, 344, 52, 0, 0, 0),
  @(181, "San Martin (Parte Francesa)", 383, 0, 273, 102, 0, 0, 8),
  @(182, "Eritrea", 375, 0, 341, 34, 0, 0, 0),
  @(183, "Curazao", 370, 0, 157, 212, 0, 0, 1),
  @(184, "Mauricio", 367, 0, 343, 14, 0, 0, 10),
  @(185, "Isla de Man", 340, 0, 315, 1, 0, 0, 24),
  @(186, "Mongolia", 313, 0, 305, 8, 0, 0, 0),
  @(187, "Butan", 281, 1, 219, 62, 0, 0, 0),
  @(188, "Camboya", 277, 0, 275, 2, 0, 0, 0),
  @(189, "Monaco", 214, 0, 178, 34, 0, 0, 2),
  @(190, "Islas Caimanes", 211, 0, 208, 2, 0, 0, 1),
  @(191, "Barbados", 190, 0, 178, 5, 0, 0, 7),
  @(192, "Bermudas", 181, 0, 167, 5, 0, 0, 9),
  @(193, "Brunei", 146, 0, 142, 1, 0, 0, 3),
  @(194, "Seychelles", 143, 0, 140, 3, 0, 0, 0),
  @(195, "Liechtenstein", 118, 0, 113, 4, 0, 0, 1),
  @(196, "Bonaire, San Eustaquio y Saba", 106, 0, 24, 81, 0, 0, 1),
  @(197, "Antigua y Barbuda", 101, 0, 92, 6, 0, 0, 3),
  @(198, "Islas Virgenes Britanicas", 71, 0, 66, 4, 0, 0, 1),
  @(199, "San Vicente y las Granadinas", 64, 0, 64, 0, 0, 0, 0),
  @(200, "San Bartolome", 48, 0, 25, 23, 0, 0, 0),
  @(201, "Macao", 46, 0, 46, 0, 0, 0, 0),
  @(202, "Puerto Rico", 39, 0, 1, 36, 0, 0, 2),
  @(203, "Guam", 32, 0, 0, 31, 0, 0, 1),
  @(204, "Fiyi", 32, 0, 28, 2, 0, 0, 2),
  @(205, "Dominica", 30, 0, 24, 6, 0, 0, 0),
  @(206, "Timor Oriental", 28, 0, 28, 0, 0, 0, 0),
  @(207, "Nueva Caledonia", 27, 0, 27, 0, 0, 0, 0),
  @(208, "Santa Lucia", 27, 0, 27, 0, 0, 0, 0),
  @(209, "Granada", 24, 0, 24, 0, 0, 0, 0),
  @(210, "Laos", 23, 0, 22, 1, 0, 0, 0),
  @(211, "San Cristobal y Nieves", 19, 0, 17, 2, 0, 0, 0),
  @(212, "Islas Virgenes de los Estados Unidos", 17, 0, 0, 17, 0, 0, 0),
  @(213, "San Pedro y Miquelon", 16, 0, 6, 10, 0, 0, 0),
  @(214, "Groenlandia", 14, 0, 14, 0, 0, 0, 0),
  @(215, "Montserrat", 13, 0, 12, 0, 0, 0, 1),
  @(216, "Islas Malvinas", 13, 0, 13, 0, 0, 0, 0),
  @(217, "Santa Sede", 12, 0, 12, 0, 0, 0, 0),
  @(218, "Sahara Occidental", 10, 0, 8, 1, 0, 0, 1),
  @(219, "Anguila", 3, 0, 3, 0, 0, 0, 0)
)

foreach ($row in $data) {
  $r = $row[0]
  $ws.Cells.Item($r, 1).Value = $row[1]
  $ws.Cells.Item($r, 2).Value = $row[2]
  $ws.Cells.Item($r, 3).Value = $row[3]
  $ws.Cells.Item($r, 4).Value = $row[4]
  $ws.Cells.Item($r, 5).Value = $row[5]
  $ws.Cells.Item($r, 6).Value = $row[6]
  $ws.Cells.Item($r, 7).Value = $row[7]
  $ws.Cells.Item($r, 8).Value = $row[8]
}

Write-Output "Done updating $($data.Count) country rows"
